$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure these cells remain text (matches source data which stores
# team names/stats as literal strings, not numbers/percentages)
$ws.Range("B2:G19").NumberFormat = "@"

$rows = @(
    @("Sporting","2.6","6.3","96%","74%","3.91"),
    @("Benfica","1.8","6.8","83%","63%","3.17"),
    @("Porto","3.0","7.3","75%","46%","2.50"),
    @("Braga","2.4","5.9","88%","71%","3.58"),
    @("Guimaraes","2.5","6.7","75%","54%","2.67"),
    @("Moreirense","2.5","4.8","50%","38%","2.21"),
    @("Arouca","2.5","3.9","79%","63%","3.17"),
    @("Gil Vicente","2.2","4.0","75%","50%","2.96"),
    @("Casa Pia","2.0","4.6","63%","38%","2.38"),
    @("Famalicão","3.4","4.8","68%","47%","2.32"),
    @("Farense","2.3","4.8","71%","54%","2.79"),
    @("Boavista","2.8","4.3","91%","61%","3.22"),
    @("Rio Ave","2.9","3.9","79%","33%","2.46"),
    @("Portimonense","3.0","3.9","79%","54%","3.04"),
    @("Estoril","2.6","5.7","83%","71%","3.58"),
    @("Estrela Amadora","2.2","4.8","71%","42%","2.54"),
    @("Chaves","3.0","4.9","88%","63%","3.25"),
    @("Vizela","2.5","4.4","67%","50%","2.92")
)

$r = 2
foreach ($row in $rows) {
    $c = 2
    foreach ($val in $row) {
        $ws.Cells.Item($r, $c).Value = $val
        $c++
    }
    $r++
}

Write-Output "Updated standings table (18 teams) for Last Update 15-03-2024"
